# Applies the "korkeakouluopinto" sample-data refresh:
#  - workbook view height tweak
#  - Koulutukset (sheet1) sample rows get new sample "slug" style values
#    instead of old free-text / example values, plus D column now mirrors C
#  - selection/active-cell bookkeeping on both sheets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Workbook-level view tweak
# ---------------------------------------------------------------------
$wb.Windows.Item(1).Height = 542

# ---------------------------------------------------------------------
# Sheet "Koulutukset"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Koulutukset")

# Column D now simply mirrors column C's OID on every data row, and loses
# its previous "text" cell style (becomes plain/general like column C).
$ws1.Range("D2:D8").Value = "1.2.246.562.10.79875033395"
$ws1.Range("D2:D8").Style = "Normal"

# Row 2 - "opintojakso a"
$ws1.Range("J2").Value = "opinnontyyppi_1"
$ws1.Range("P2").Value = "kieli_fi, kieli_sv, kieli_en"
$ws1.Range("R2").Value = "opetusaikakk_2"
$ws1.Range("S2").Value = "opetusmuoto_e"
$ws1.Range("T2").Value = "opetuspaikkakk_2"
$ws1.Range("U2").Value = "teemat_1"
$ws1.Range("V2").Value = "oppiaineetyleissivistava"

# Row 3 - "opintojakso b"
$ws1.Range("J3").Value = "opinnontyyppi_1"
$ws1.Range("P3").Value = "kieli_sv"
$ws1.Range("R3").Value = "opetusaikakk_1"
$ws1.Range("S3").Value = "opetusmuoto_e"
$ws1.Range("T3").Value = "opetuspaikkakk_2"
$ws1.Range("U3").Value = "teemat_1, teemat_2"

# Row 4 - "opintojakso c"
$ws1.Range("J4").Value = "opinnontyyppi_1"
$ws1.Range("P4").Value = "kieli_fi, kieli_sv, kieli_en"
$ws1.Range("R4").Value = "opetusaikakk_1"
$ws1.Range("S4").Value = "opetusmuoto_e"
$ws1.Range("T4").Value = "opetuspaikkakk_2"

# Row 5 - "opintokokonaisuus abc"
$ws1.Range("J5").Value = ""
$ws1.Range("P5").Value = ""
$ws1.Range("R5").Value = "opetusaikakk_1"
$ws1.Range("S5").Value = "opetusmuoto_e"
$ws1.Range("T5").Value = "opetuspaikkakk_2"

# Row 6 - "opintojakso x"
$ws1.Range("P6").Value = ""
$ws1.Range("R6").Value = "opetusaikakk_1"
$ws1.Range("S6").Value = "opetusmuoto_e"
$ws1.Range("T6").Value = "opetuspaikkakk_2"

# Row 7 - "opintokokonaisuus xy"
$ws1.Range("J7").Value = ""

# Row 8 - "opintokokonaisuus xy"
$ws1.Range("J8").Value = ""

# Active selection moved to B8 on this sheet
$ws1.Range("B8").Select()

# ---------------------------------------------------------------------
# Sheet "Hakukohteet"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Hakukohteet")

# Active selection moved to D3 on this sheet
$ws2.Range("D3").Select()

$wb.Save()
